$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.282.11"
$ws.Range("E2").Value = "  -2.58%  "
$ws.Range("D3").Value = "1.567.77"
$ws.Range("E3").Value = "  -3.45%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.68"
$ws.Range("E5").Value = "  -2.76%  "
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.477"
$ws.Range("E7").Value = "  -5.04%  "
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0607"
$ws.Range("E8").Value = "  -1.74%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.243"
$ws.Range("E9").Value = "  -2.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.83"
$ws.Range("E10").Value = "  -2.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0781"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("D12").Value = "1.786.99"
$ws.Range("E12").Value = "  -3.31%  "
$ws.Range("D13").Value = "1.584.04"
$ws.Range("E13").Value = "  -2.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.02"
$ws.Range("E14").Value = "  -3.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.506"
$ws.Range("E15").Value = "  -3.43%  "
$ws.Range("D16").Value = "25.301.95"
$ws.Range("E16").Value = "  -2.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "59.46"
$ws.Range("E17").Value = "  -2.80%  "
$ws.Range("E18").Value = "  -3.03%  "
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "185.88"
$ws.Range("E20").Value = "  -2.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.14"
$ws.Range("E21").Value = "  -2.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.29"
$ws.Range("E22").Value = "  -2.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.87"
$ws.Range("E23").Value = "  -3.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.129"
$ws.Range("E24").Value = "  -3.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.01"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.54"
$ws.Range("E26").Value = "  -2.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.65"
$ws.Range("E27").Value = "  -5.95%  "
$ws.Range("E28").Value = "  -3.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "14.86"
$ws.Range("E29").Value = "  -2.08%  "
$ws.Range("E30").Value = "  -5.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0464"
$ws.Range("E31").Value = "  -3.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.04"
$ws.Range("E32").Value = "  -2.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.00"
$ws.Range("E33").Value = "  -3.59%  "
$ws.Range("E34").Value = "  -1.93%  "
$ws.Range("E35").Value = "  -4.00%  "
$ws.Range("D36").Value = "1.090.45"
$ws.Range("E36").Value = "  -3.08%  "
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("E38").Value = "  -4.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0150"
$ws.Range("E39").Value = "  -2.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.496"
$ws.Range("E40").Value = "  -3.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.775"
$ws.Range("E41").Value = "  -8.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.765"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "93.36"
$ws.Range("E43").Value = "  -4.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.06"
$ws.Range("E44").Value = "  -1.97%  "
$ws.Range("D45").Value = "1.701.25"
$ws.Range("E45").Value = "  -3.29%  "
$ws.Range("D46").Value = "0.0₆0111"
$ws.Range("E46").Value = "  -2.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "52.80"
$ws.Range("E47").Value = "  -3.01%  "
$ws.Range("E48").Value = "  -4.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.43"
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.407"
$ws.Range("E50").Value = "  -1.59%  "
$ws.Range("E51").Value = "  -0.51%  "
